# Auto-generated edit script applying the cryptos.xlsx diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must stay TEXT even when it parses as a number.
# Excel auto-converts plain numeric-looking strings ("5.42", "0.998", ...) to
# real numbers on assignment. Prefixing the literal with an apostrophe forces
# text entry (like typing '5.42 into a cell). That marks the cell with a
# "quote prefix" style, so we reset Style back to "Normal" afterwards to
# leave the cell on the default/general style, matching a plain text cell.
function Set-TextValue([string]$rangeAddr, [string]$text) {
    $ws.Range($rangeAddr).Value = "'" + $text
    $ws.Range($rangeAddr).Style = "Normal"
}

$ws.Range("D2").Value = '61.243.31'
$ws.Range("E2").Value = '  -4.24%  '
$ws.Range("D3").Value = '2.453.52'
$ws.Range("E3").Value = '  -7.05%  '
$ws.Range("E4").Value = '  -0.03%  '
Set-TextValue "D5" '549.46'
$ws.Range("E5").Value = '  -5.13%  '
Set-TextValue "D6" '146.09'
$ws.Range("E6").Value = '  -6.95%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  -6.22%  '
$ws.Range("D9").Value = '2.448.87'
$ws.Range("E9").Value = '  -7.18%  '
$ws.Range("E10").Value = '  -10.05%  '
$ws.Range("B11").Value = 'TRON'
$ws.Range("C11").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue "D11" '0.154'
$ws.Range("E11").Value = '  -1.86%  '
$ws.Range("B12").Value = 'Toncoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue "D12" '5.42'
$ws.Range("E12").Value = '  -7.21%  '
$ws.Range("E13").Value = '  -8.34%  '
Set-TextValue "D14" '26.01'
$ws.Range("E14").Value = '  -9.34%  '
$ws.Range("D15").Value = '2.895.59'
$ws.Range("E15").Value = '  -7.04%  '
$ws.Range("E16").Value = '  -10.68%  '
$ws.Range("D17").Value = '61.141.57'
$ws.Range("E17").Value = '  -4.25%  '
$ws.Range("D18").Value = '2.457.00'
$ws.Range("E18").Value = '  -7.41%  '
$ws.Range("E19").Value = '  -9.16%  '
Set-TextValue "D20" '7.15'
$ws.Range("E20").Value = '  -7.68%  '
$ws.Range("E21").Value = '  -7.92%  '
Set-TextValue "D22" '317.50'
$ws.Range("E22").Value = '  -7.95%  '
$ws.Range("E23").Value = '  +0.00%  '
$ws.Range("E24").Value = '  -0.59%  '
Set-TextValue "D25" '63.88'
$ws.Range("E25").Value = '  -6.47%  '
$ws.Range("E26").Value = '  -13.22%  '
$ws.Range("D27").Value = '2.598.11'
$ws.Range("E27").Value = '  -6.21%  '
Set-TextValue "D28" '551.09'
$ws.Range("E28").Value = '  -5.28%  '
$ws.Range("E29").Value = '  +0.13%  '
$ws.Range("E30").Value = '  -10.22%  '
$ws.Range("E31").Value = '  -11.24%  '
Set-TextValue "D32" '7.62'
$ws.Range("E32").Value = '  -7.48%  '
$ws.Range("E33").Value = '  -8.97%  '
$ws.Range("E34").Value = '  -7.64%  '
$ws.Range("E35").Value = '  -9.05%  '
$ws.Range("E36").Value = '  -11.66%  '
Set-TextValue "D37" '0.998'
$ws.Range("E37").Value = '  -0.09%  '
Set-TextValue "D38" '4.82'
$ws.Range("E38").Value = '  -11.72%  '
Set-TextValue "D39" '0.378'
$ws.Range("E39").Value = '  -6.32%  '
Set-TextValue "D40" '18.34'
$ws.Range("E40").Value = '  -7.32%  '
Set-TextValue "D41" '141.56'
$ws.Range("E41").Value = '  -7.36%  '
Set-TextValue "D42" '1.74'
$ws.Range("E42").Value = '  -9.10%  '
$ws.Range("E43").Value = '  +0.07%  '
Set-TextValue "D44" '40.40'
$ws.Range("E44").Value = '  -4.33%  '
Set-TextValue "D45" '2.39'
$ws.Range("E45").Value = '  -6.46%  '
Set-TextValue "D46" '146.29'
$ws.Range("E46").Value = '  -10.02%  '
$ws.Range("E47").Value = '  -8.42%  '
Set-TextValue "D48" '21.33'
$ws.Range("E48").Value = '  -12.14%  '
$ws.Range("E49").Value = '  -9.59%  '
Set-TextValue "D50" '0.588'
$ws.Range("E50").Value = '  -7.39%  '
$ws.Range("E51").Value = '  -7.06%  '
